$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the original data range so that every value is rewritten fresh in row-major
# order below; this keeps the shared-string table build order aligned with the target.
$ws.Range("A2:K12").ClearContents()

# Write all data rows (TestScenario_1..4) fresh, in row-major order
$ws.Range("A2").Value = "TestScenario_1"
$ws.Range("B2").Value = "TestScenario_1.TestCase_1"
$ws.Range("C2").Value = "New Case"
$ws.Range("D2").Value = "User Needs to Login to Salesforce, from the browser with correct credentials"
$ws.Range("F2").Value = "Step 1"
$ws.Range("G2").Value = "Click on the Case tab, and click on New button"
$ws.Range("H2").Value = "User should be navigated to the New  Case Page"
$ws.Range("I2").Value = "Approved"
$ws.Range("E3").Value = "Valid value for required field Direction "
$ws.Range("F3").Value = "Step 2"
$ws.Range("G3").Value = "Input valid value in the  Direction field."
$ws.Range("H3").Value = "User should be able to input value for the Direction field."
$ws.Range("E4").Value = "Valid value for required field Status, value should be equals New to submit for Approval towards to the assigned approver  @ Valid value for required field Status, value should be  Responded for Workflow Process."
$ws.Range("F4").Value = "Step 3"
$ws.Range("G4").Value = "Input valid value in the  Status field."
$ws.Range("H4").Value = "Value accepted for Status field."
$ws.Range("E5").Value = "Valid value for required field Account Id "
$ws.Range("F5").Value = "Step 4"
$ws.Range("G5").Value = "Input valid value in the  Account Id field."
$ws.Range("H5").Value = "User should be able to input value for the Account Id field."
$ws.Range("E6").Value = "Valid value for required field Root Cause "
$ws.Range("F6").Value = "Step 5"
$ws.Range("G6").Value = "Input valid value in the  Root Cause field."
$ws.Range("H6").Value = "User should be able to input value for the Root Cause field."
$ws.Range("F7").Value = "Step 6"
$ws.Range("G7").Value = "Click on Save button to save Case with fields"
$ws.Range("H7").Value = "User should be able to validate that a New Case is created"
$ws.Range("A8").Value = "TestScenario_2"
$ws.Range("B8").Value = "TestScenario_2.TestCase_1"
$ws.Range("C8").Value = "View Case"
$ws.Range("D8").Value = "User Needs to Login to Salesforce, from the browser with correct credentials"
$ws.Range("F8").Value = "Step 1"
$ws.Range("G8").Value = "Click on the Case tab,  and select a Case "
$ws.Range("H8").Value = "User should be navigated to the Case Page"
$ws.Range("I8").Value = "Approved"
$ws.Range("F9").Value = "Step 2"
$ws.Range("G9").Value = "Click on Case name to View the Details"
$ws.Range("H9").Value = "User should be able to view the Case Details"
$ws.Range("A10").Value = "TestScenario_3"
$ws.Range("B10").Value = "TestScenario_3.TestCase_1"
$ws.Range("C10").Value = "Edit Case"
$ws.Range("D10").Value = "User Needs to Login to Salesforce, from the browser with correct credentials"
$ws.Range("F10").Value = "Step 1"
$ws.Range("G10").Value = "Click on the Case tab,  and click on existing Case to modify"
$ws.Range("H10").Value = "User is navigated to the Case Details page"
$ws.Range("I10").Value = "Approved"
$ws.Range("E11").Value = "Valid value for required field Direction "
$ws.Range("F11").Value = "Step 2"
$ws.Range("G11").Value = "Input valid value in the  Direction field."
$ws.Range("H11").Value = "User should be able to input value for the Direction field."
$ws.Range("E12").Value = "Valid value for required field Status, value should be equals New to submit for Approval towards to the assigned approver  @ Valid value for required field Status, value should be  Responded for Workflow Process."
$ws.Range("F12").Value = "Step 3"
$ws.Range("G12").Value = "Input valid value in the  Status field."
$ws.Range("H12").Value = "Value accepted for Status field."
$ws.Range("E13").Value = "Valid value for required field Account Id "
$ws.Range("F13").Value = "Step 4"
$ws.Range("G13").Value = "Input valid value in the  Account Id field."
$ws.Range("H13").Value = "User should be able to input value for the Account Id field."
$ws.Range("E14").Value = "Valid value for required field Root Cause "
$ws.Range("F14").Value = "Step 5"
$ws.Range("G14").Value = "Input valid value in the  Root Cause field."
$ws.Range("H14").Value = "User should be able to input value for the Root Cause field."
$ws.Range("F15").Value = "Step 6"
$ws.Range("G15").Value = "Click on Save button to save Case with fields"
$ws.Range("H15").Value = "User should be able to validate that the Case is edited"
$ws.Range("A16").Value = "TestScenario_4"
$ws.Range("B16").Value = "TestScenario_4.TestCase_1"
$ws.Range("C16").Value = "Delete Case"
$ws.Range("D16").Value = "User Needs to Login to Salesforce, from the browser with correct credentials"
$ws.Range("F16").Value = "Step 1"
$ws.Range("G16").Value = "Click on the Case tab,  and select the existing  Case to delete"
$ws.Range("H16").Value = "User is navigated to the Case Details page"
$ws.Range("I16").Value = "Approved"
$ws.Range("F17").Value = "Step 2"
$ws.Range("G17").Value = "Click on to the Delete to Delete the Case"
$ws.Range("H17").Value = "User should be able to validate that a pop-up is displayed asking for confirmation to delete the Case"
$ws.Range("F18").Value = "Step 3"
$ws.Range("G18").Value = "Click on Confirm / OK to delete the  Case"
$ws.Range("H18").Value = "User should be able to validate the Case is deleted"

# Resize the table (ListObject) to cover the new row range
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:K18"))

# Adjust column widths to match the updated layout
$ws.Columns.Item(3).ColumnWidth = 12.166666666666666
$ws.Columns.Item(5).ColumnWidth = 189.0
$ws.Columns.Item(7).ColumnWidth = 55.0
$ws.Columns.Item(8).ColumnWidth = 89.66666666666667
